$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2021", 2)

# 2. Update docente responsável
$d.Content.Find.Execute("5840560 - Marco Antonio Carvalho Pereira", $true, $false, $false, $false, $false,
                         $true, 1, $false, "11079086 - Herlandí de Souza Andrade", 2)

# 3. Update critério de avaliação
$d.Content.Find.Execute("Média aritmética de duas provas teóricas.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Média Aritmética das atividades avaliativas realizadas.", 2)

# 4. Remove trailing period from norma de recuperação
$d.Content.Find.Execute("Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação", 2)

# 5. Add a new requisito fraco line after the LOQ4205 line, inside the same paragraph
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$insertPoint.InsertAfter("LOQ4240 -  Administração e Organização II  (Requisito fraco)" + [char]11)
